$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first two data records (rows 2 and 3) are gone. Column A keeps its
# running 0-based index as-is, but the "codice_particella" / "codice_comune
# _catastale" values in columns B:C need to shift up by two rows to close
# the gap, so copy B4:C73 over B2:C3... and onward, then drop the two
# leftover rows at the bottom.
$ws.Range("B4:C73").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Rows("72:73").Delete()
